$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New voltage magnitude values per bus (case with 380 kV)
# Row 2
$ws.Cells.Item(2, 2).Value2 = 1.02
$ws.Cells.Item(2, 3).Value2 = 1.043932743268985
$ws.Cells.Item(2, 4).Value2 = 1.045237337436435
$ws.Cells.Item(2, 5).Value2 = 1.041800697123538
$ws.Cells.Item(2, 6).Value2 = 1.049099777258221
$ws.Cells.Item(2, 9).Value2 = 1.044566932928541
$ws.Cells.Item(2, 10).Value2 = 1.049000537369658
$ws.Cells.Item(2, 11).Value2 = 1.048006141244341
$ws.Cells.Item(2, 12).Value2 = 1.044579186282664
$ws.Cells.Item(2, 13).Value2 = 1.051857777827612
$ws.Cells.Item(2, 14).Value2 = 1.050490238454766

# Row 3
$ws.Cells.Item(3, 2).Value2 = 1.02
$ws.Cells.Item(3, 3).Value2 = 1.045143125510318
$ws.Cells.Item(3, 4).Value2 = 1.04617860050535
$ws.Cells.Item(3, 5).Value2 = 1.042844677160427
$ws.Cells.Item(3, 6).Value2 = 1.051487736312702
$ws.Cells.Item(3, 9).Value2 = 1.045044182254169
$ws.Cells.Item(3, 10).Value2 = 1.04985667644488
$ws.Cells.Item(3, 11).Value2 = 1.048758465866926
$ws.Cells.Item(3, 12).Value2 = 1.045433256027789
$ws.Cells.Item(3, 13).Value2 = 1.054053848853092
$ws.Cells.Item(3, 14).Value2 = 1.051347593345676

# Row 4
$ws.Cells.Item(4, 2).Value2 = 1.02
$ws.Cells.Item(4, 3).Value2 = 1.045925087282464
$ws.Cells.Item(4, 4).Value2 = 1.046786432886435
$ws.Cells.Item(4, 5).Value2 = 1.043519237544882
$ws.Cells.Item(4, 6).Value2 = 1.05302650123891
$ws.Cells.Item(4, 9).Value2 = 1.045350669819362
$ws.Cells.Item(4, 10).Value2 = 1.050408909171509
$ws.Cells.Item(4, 11).Value2 = 1.049243422393772
$ws.Cells.Item(4, 12).Value2 = 1.04598433441069
$ws.Cells.Item(4, 13).Value2 = 1.055468156743401
$ws.Cells.Item(4, 14).Value2 = 1.051900610306119

# Row 5
$ws.Cells.Item(5, 2).Value2 = 1.02
$ws.Cells.Item(5, 3).Value2 = 1.046253533302286
$ws.Cells.Item(5, 4).Value2 = 1.047041675676442
$ws.Cells.Item(5, 5).Value2 = 1.043802595812555
$ws.Cells.Item(5, 6).Value2 = 1.053671903823833
$ws.Cells.Item(5, 9).Value2 = 1.045478964563315
$ws.Cells.Item(5, 10).Value2 = 1.050640653908222
$ws.Cells.Item(5, 11).Value2 = 1.049446859704908
$ws.Cells.Item(5, 12).Value2 = 1.046215637707377
$ws.Cells.Item(5, 13).Value2 = 1.056061165033712
$ws.Cells.Item(5, 14).Value2 = 1.052132684146939

# Row 6
$ws.Cells.Item(6, 2).Value2 = 1.02
$ws.Cells.Item(6, 3).Value2 = 1.046308663939979
$ws.Cells.Item(6, 4).Value2 = 1.047084515176395
$ws.Cells.Item(6, 5).Value2 = 1.043850159720871
$ws.Cells.Item(6, 6).Value2 = 1.0537801832605
$ws.Cells.Item(6, 9).Value2 = 1.045500473521251
$ws.Cells.Item(6, 10).Value2 = 1.050679540719315
$ws.Cells.Item(6, 11).Value2 = 1.049480992129858
$ws.Cells.Item(6, 12).Value2 = 1.046254452959555
$ws.Cells.Item(6, 13).Value2 = 1.05616064291109
$ws.Cells.Item(6, 14).Value2 = 1.052171626181764

# Row 7
$ws.Cells.Item(7, 2).Value2 = 1.02
$ws.Cells.Item(7, 3).Value2 = 1.045929477129026
$ws.Cells.Item(7, 4).Value2 = 1.046789844586696
$ws.Cells.Item(7, 5).Value2 = 1.043523024677851
$ws.Cells.Item(7, 6).Value2 = 1.053035130958695
$ws.Cells.Item(7, 9).Value2 = 1.045352386266048
$ws.Cells.Item(7, 10).Value2 = 1.050412007376413
$ws.Cells.Item(7, 11).Value2 = 1.04924614245016
$ws.Cells.Item(7, 12).Value2 = 1.045987426545245
$ws.Cells.Item(7, 13).Value2 = 1.055476086656485
$ws.Cells.Item(7, 14).Value2 = 1.051903712910829

# Row 8
$ws.Cells.Item(8, 2).Value2 = 1.02
$ws.Cells.Item(8, 3).Value2 = 1.04434205684541
$ws.Cells.Item(8, 4).Value2 = 1.045555698330704
$ws.Cells.Item(8, 5).Value2 = 1.042153717188038
$ws.Cells.Item(8, 6).Value2 = 1.049908148882319
$ws.Cells.Item(8, 9).Value2 = 1.044728705272684
$ws.Cells.Item(8, 10).Value2 = 1.049290237756808
$ws.Cells.Item(8, 11).Value2 = 1.048260777871785
$ws.Cells.Item(8, 12).Value2 = 1.044868148878568
$ws.Cells.Item(8, 13).Value2 = 1.05260135907733
$ws.Cells.Item(8, 14).Value2 = 1.050780350249695

# Row 9
$ws.Cells.Item(9, 2).Value2 = 1.02
$ws.Cells.Item(9, 3).Value2 = 1.041535109370893
$ws.Cells.Item(9, 4).Value2 = 1.043371390555874
$ws.Cells.Item(9, 5).Value2 = 1.039733254481126
$ws.Cells.Item(9, 6).Value2 = 1.044347189931047
$ws.Cells.Item(9, 9).Value2 = 1.043611718513763
$ws.Cells.Item(9, 10).Value2 = 1.047299969002032
$ws.Cells.Item(9, 11).Value2 = 1.046510107938976
$ws.Cells.Item(9, 12).Value2 = 1.04288370180426
$ws.Cells.Item(9, 13).Value2 = 1.047482776394319
$ws.Cells.Item(9, 14).Value2 = 1.048787255085001

# Row 10
$ws.Cells.Item(10, 2).Value2 = 1.02
$ws.Cells.Item(10, 3).Value2 = 1.039656941996355
$ws.Cells.Item(10, 4).Value2 = 1.041908483326912
$ws.Cells.Item(10, 5).Value2 = 1.03811426685003
$ws.Cells.Item(10, 6).Value2 = 1.040603222881471
$ws.Cells.Item(10, 9).Value2 = 1.042854720035505
$ws.Cells.Item(10, 10).Value2 = 1.045963729567423
$ws.Cells.Item(10, 11).Value2 = 1.045333103565895
$ws.Cells.Item(10, 12).Value2 = 1.041552328367055
$ws.Cells.Item(10, 13).Value2 = 1.044032454904197
$ws.Cells.Item(10, 14).Value2 = 1.047449118037126

# Row 11
$ws.Cells.Item(11, 2).Value2 = 1.02
$ws.Cells.Item(11, 3).Value2 = 1.038841968357594
$ws.Cells.Item(11, 4).Value2 = 1.041273379865988
$ws.Cells.Item(11, 5).Value2 = 1.037411900695516
$ws.Cells.Item(11, 6).Value2 = 1.038972786639726
$ws.Cells.Item(11, 9).Value2 = 1.042523947338278
$ws.Cells.Item(11, 10).Value2 = 1.045382833353632
$ws.Cells.Item(11, 11).Value2 = 1.044821043585509
$ws.Cells.Item(11, 12).Value2 = 1.040973777979647
$ws.Cells.Item(11, 13).Value2 = 1.042528904453355
$ws.Cells.Item(11, 14).Value2 = 1.046867396884082

# Row 12
$ws.Cells.Item(12, 2).Value2 = 1.02
$ws.Cells.Item(12, 3).Value2 = 1.038538986451675
$ws.Cells.Item(12, 4).Value2 = 1.041037220318412
$ws.Cells.Item(12, 5).Value2 = 1.037150805055709
$ws.Cells.Item(12, 6).Value2 = 1.038365727961106
$ws.Cells.Item(12, 9).Value2 = 1.042400629798475
$ws.Cells.Item(12, 10).Value2 = 1.045166712281658
$ws.Cells.Item(12, 11).Value2 = 1.044630474611011
$ws.Cells.Item(12, 12).Value2 = 1.040758564584194
$ws.Cells.Item(12, 13).Value2 = 1.041968940128231
$ws.Cells.Item(12, 14).Value2 = 1.046650968895395

# Row 13
$ws.Cells.Item(13, 2).Value2 = 1.02
$ws.Cells.Item(13, 3).Value2 = 1.038603989174042
$ws.Cells.Item(13, 4).Value2 = 1.041087888914479
$ws.Cells.Item(13, 5).Value2 = 1.037206820343668
$ws.Cells.Item(13, 6).Value2 = 1.038496010157282
$ws.Cells.Item(13, 9).Value2 = 1.042427102447661
$ws.Cells.Item(13, 10).Value2 = 1.045213086927856
$ws.Cells.Item(13, 11).Value2 = 1.044671368998086
$ws.Cells.Item(13, 12).Value2 = 1.040804742879238
$ws.Cells.Item(13, 13).Value2 = 1.042089122035381
$ws.Cells.Item(13, 14).Value2 = 1.04669740939891

# Row 14
$ws.Cells.Item(14, 2).Value2 = 1.02
$ws.Cells.Item(14, 3).Value2 = 1.038816929212205
$ws.Cells.Item(14, 4).Value2 = 1.041253864052873
$ws.Cells.Item(14, 5).Value2 = 1.037390322672478
$ws.Cells.Item(14, 6).Value2 = 1.038922636683057
$ws.Cells.Item(14, 9).Value2 = 1.042513763159331
$ws.Cells.Item(14, 10).Value2 = 1.045364975897336
$ws.Cells.Item(14, 11).Value2 = 1.044805298613272
$ws.Cells.Item(14, 12).Value2 = 1.040955994808362
$ws.Cells.Item(14, 13).Value2 = 1.042482648060822
$ws.Cells.Item(14, 14).Value2 = 1.046849514068149

# Row 15
$ws.Cells.Item(15, 2).Value2 = 1.02
$ws.Cells.Item(15, 3).Value2 = 1.038948093356147
$ws.Cells.Item(15, 4).Value2 = 1.041356093019529
$ws.Cells.Item(15, 5).Value2 = 1.03750335711612
$ws.Cells.Item(15, 6).Value2 = 1.039185302803944
$ws.Cells.Item(15, 9).Value2 = 1.042567097395259
$ws.Cells.Item(15, 10).Value2 = 1.04545851311468
$ws.Cells.Item(15, 11).Value2 = 1.04488776827057
$ws.Cells.Item(15, 12).Value2 = 1.041049144352163
$ws.Cells.Item(15, 13).Value2 = 1.042724915103968
$ws.Cells.Item(15, 14).Value2 = 1.04694318411907

# Row 16
$ws.Cells.Item(16, 2).Value2 = 1.02
$ws.Cells.Item(16, 3).Value2 = 1.039710992407511
$ws.Cells.Item(16, 4).Value2 = 1.041950597741423
$ws.Cells.Item(16, 5).Value2 = 1.03816085203563
$ws.Cells.Item(16, 6).Value2 = 1.040711230033435
$ws.Cells.Item(16, 9).Value2 = 1.042876608963642
$ws.Cells.Item(16, 10).Value2 = 1.046002232904762
$ws.Cells.Item(16, 11).Value2 = 1.045367036103199
$ws.Cells.Item(16, 12).Value2 = 1.041590681089668
$ws.Cells.Item(16, 13).Value2 = 1.044132035749818
$ws.Cells.Item(16, 14).Value2 = 1.04748767605362

# Row 17
$ws.Cells.Item(17, 2).Value2 = 1.02
$ws.Cells.Item(17, 3).Value2 = 1.040189074973631
$ws.Cells.Item(17, 4).Value2 = 1.042323068116899
$ws.Cells.Item(17, 5).Value2 = 1.038572920489423
$ws.Cells.Item(17, 6).Value2 = 1.041665886852687
$ws.Cells.Item(17, 9).Value2 = 1.043069954077834
$ws.Cells.Item(17, 10).Value2 = 1.046342675515206
$ws.Cells.Item(17, 11).Value2 = 1.045667019516814
$ws.Cells.Item(17, 12).Value2 = 1.041929818583974
$ws.Cells.Item(17, 13).Value2 = 1.045012099637897
$ws.Cells.Item(17, 14).Value2 = 1.04782860213162

# Row 18
$ws.Cells.Item(18, 2).Value2 = 1.02
$ws.Cells.Item(18, 3).Value2 = 1.040467767459407
$ws.Cells.Item(18, 4).Value2 = 1.042540164460704
$ws.Cells.Item(18, 5).Value2 = 1.038813144707244
$ws.Cells.Item(18, 6).Value2 = 1.042221830574884
$ws.Cells.Item(18, 9).Value2 = 1.043182441118282
$ws.Cells.Item(18, 10).Value2 = 1.046541028799362
$ws.Cells.Item(18, 11).Value2 = 1.045841762653357
$ws.Cells.Item(18, 12).Value2 = 1.042127433627415
$ws.Cells.Item(18, 13).Value2 = 1.045524508817067
$ws.Cells.Item(18, 14).Value2 = 1.048027237100191

# Row 19
$ws.Cells.Item(19, 2).Value2 = 1.02
$ws.Cells.Item(19, 3).Value2 = 1.040562766613909
$ws.Cells.Item(19, 4).Value2 = 1.042614161835233
$ws.Cells.Item(19, 5).Value2 = 1.038895033344241
$ws.Cells.Item(19, 6).Value2 = 1.042411243275232
$ws.Cells.Item(19, 9).Value2 = 1.043220747610737
$ws.Cells.Item(19, 10).Value2 = 1.046608624799006
$ws.Cells.Item(19, 11).Value2 = 1.04590130633027
$ws.Cells.Item(19, 12).Value2 = 1.042194781795617
$ws.Cells.Item(19, 13).Value2 = 1.045699072903332
$ws.Cells.Item(19, 14).Value2 = 1.048094929093908

# Row 20
$ws.Cells.Item(20, 2).Value2 = 1.02
$ws.Cells.Item(20, 3).Value2 = 1.040137798353203
$ws.Cells.Item(20, 4).Value2 = 1.042283122087898
$ws.Cells.Item(20, 5).Value2 = 1.038528722756372
$ws.Cells.Item(20, 6).Value2 = 1.04156355376789
$ws.Cells.Item(20, 9).Value2 = 1.043049239791206
$ws.Cells.Item(20, 10).Value2 = 1.046306172156957
$ws.Cells.Item(20, 11).Value2 = 1.045634858181119
$ws.Cells.Item(20, 12).Value2 = 1.041893452868905
$ws.Cells.Item(20, 13).Value2 = 1.044917772309785
$ws.Cells.Item(20, 14).Value2 = 1.047792046934415

# Row 21
$ws.Cells.Item(21, 2).Value2 = 1.02
$ws.Cells.Item(21, 3).Value2 = 1.038754231032325
$ws.Cells.Item(21, 4).Value2 = 1.041204995561
$ws.Cells.Item(21, 5).Value2 = 1.037336291527297
$ws.Cells.Item(21, 6).Value2 = 1.038797046008756
$ws.Cells.Item(21, 9).Value2 = 1.042488256308398
$ws.Cells.Item(21, 10).Value2 = 1.045320258097211
$ws.Cells.Item(21, 11).Value2 = 1.044765869854992
$ws.Cells.Item(21, 12).Value2 = 1.040911463594124
$ws.Cells.Item(21, 13).Value2 = 1.04236680564632
$ws.Cells.Item(21, 14).Value2 = 1.046804732763618

# Row 22
$ws.Cells.Item(22, 2).Value2 = 1.02
$ws.Cells.Item(22, 3).Value2 = 1.037882794119753
$ws.Cells.Item(22, 4).Value2 = 1.040525663316836
$ws.Cells.Item(22, 5).Value2 = 1.036585370652073
$ws.Cells.Item(22, 6).Value2 = 1.037049264629452
$ws.Cells.Item(22, 9).Value2 = 1.042132915446745
$ws.Cells.Item(22, 10).Value2 = 1.044698344996942
$ws.Cells.Item(22, 11).Value2 = 1.044217375886774
$ws.Cells.Item(22, 12).Value2 = 1.040292228598785
$ws.Cells.Item(22, 13).Value2 = 1.040754332865615
$ws.Cells.Item(22, 14).Value2 = 1.046181936475413

# Row 23
$ws.Cells.Item(23, 2).Value2 = 1.02
$ws.Cells.Item(23, 3).Value2 = 1.038344906936859
$ws.Cells.Item(23, 4).Value2 = 1.040885931415271
$ws.Cells.Item(23, 5).Value2 = 1.036983562758817
$ws.Cells.Item(23, 6).Value2 = 1.037976606331729
$ws.Cells.Item(23, 9).Value2 = 1.042321539187814
$ws.Cells.Item(23, 10).Value2 = 1.045028227140776
$ws.Cells.Item(23, 11).Value2 = 1.044508346219099
$ws.Cells.Item(23, 12).Value2 = 1.040620670928277
$ws.Cells.Item(23, 13).Value2 = 1.04160996384129
$ws.Cells.Item(23, 14).Value2 = 1.046512287089729

# Row 24
$ws.Cells.Item(24, 2).Value2 = 1.02
$ws.Cells.Item(24, 3).Value2 = 1.040160968560829
$ws.Cells.Item(24, 4).Value2 = 1.042301172471449
$ws.Cells.Item(24, 5).Value2 = 1.038548694206278
$ws.Cells.Item(24, 6).Value2 = 1.041609796436498
$ws.Cells.Item(24, 9).Value2 = 1.043058600575264
$ws.Cells.Item(24, 10).Value2 = 1.046322667136715
$ws.Cells.Item(24, 11).Value2 = 1.04564939122169
$ws.Cells.Item(24, 12).Value2 = 1.041909885582557
$ws.Cells.Item(24, 13).Value2 = 1.044960397601248
$ws.Cells.Item(24, 14).Value2 = 1.047808565338936

# Row 25
$ws.Cells.Item(25, 2).Value2 = 1.02
$ws.Cells.Item(25, 3).Value2 = 1.042261958045016
$ws.Cells.Item(25, 4).Value2 = 1.043937248571208
$ws.Cells.Item(25, 5).Value2 = 1.040359925221447
$ws.Cells.Item(25, 6).Value2 = 1.045791107396528
$ws.Cells.Item(25, 9).Value2 = 1.043902642823013
$ws.Cells.Item(25, 10).Value2 = 1.047816136693884
$ws.Cells.Item(25, 11).Value2 = 1.046964422484468
$ws.Cells.Item(25, 12).Value2 = 1.043398192576815
$ws.Cells.Item(25, 13).Value2 = 1.049304155794179
$ws.Cells.Item(25, 14).Value2 = 1.049304155794179
